$wb = $excel.ActiveWorkbook

# --- sheet index 1 (COM item 2) : new column 73 ---
$ws = $wb.Worksheets.Item(2)
$ws.Columns.Item(73).ColumnWidth = $ws.Columns.Item(72).ColumnWidth
$ws.Cells.Item(1, 73).Value = 20251209
$ws.Cells.Item(1, 72).Copy()
$ws.Cells.Item(1, 73).PasteSpecial(-4122)
$ws.Cells.Item(2, 73).Value = 278.16
$ws.Cells.Item(3, 73).Value = 226.84
$ws.Cells.Item(4, 73).Value = 312.37
$ws.Cells.Item(5, 73).Value = 489.1
$ws.Cells.Item(6, 73).Value = 97.03
$ws.Cells.Item(7, 73).Value = 185.56
$ws.Cells.Item(8, 73).Value = 180.97
$ws.Cells.Item(9, 73).Value = 83.32
$ws.Cells.Item(10, 73).Value = 437.54

# --- sheet index 2 (COM item 3) : new column 73 ---
$ws = $wb.Worksheets.Item(3)
$ws.Columns.Item(73).ColumnWidth = $ws.Columns.Item(72).ColumnWidth
$ws.Cells.Item(1, 73).Value = 20251209
$ws.Cells.Item(1, 72).Copy()
$ws.Cells.Item(1, 73).PasteSpecial(-4122)
$ws.Cells.Item(2, 73).Value = 280.03
$ws.Cells.Item(3, 73).Value = 228.57
$ws.Cells.Item(4, 73).Value = 317.99
$ws.Cells.Item(5, 73).Value = 492.12
$ws.Cells.Item(6, 73).Value = 97.24
$ws.Cells.Item(7, 73).Value = 185.7167
$ws.Cells.Item(8, 73).Value = 182.8499
$ws.Cells.Item(9, 73).Value = 83.33
$ws.Cells.Item(10, 73).Value = 452.39

# --- sheet index 3 (COM item 4) : new column 73 ---
$ws = $wb.Worksheets.Item(4)
$ws.Columns.Item(73).ColumnWidth = $ws.Columns.Item(72).ColumnWidth
$ws.Cells.Item(1, 73).Value = 20251209
$ws.Cells.Item(1, 72).Copy()
$ws.Cells.Item(1, 73).PasteSpecial(-4122)
$ws.Cells.Item(2, 73).Value = 276.92
$ws.Cells.Item(3, 73).Value = 225.11
$ws.Cells.Item(4, 73).Value = 311.895
$ws.Cells.Item(5, 73).Value = 488.5001
$ws.Cells.Item(6, 73).Value = 95.45
$ws.Cells.Item(7, 73).Value = 183.32
$ws.Cells.Item(8, 73).Value = 180.38
$ws.Cells.Item(9, 73).Value = 82.02
$ws.Cells.Item(10, 73).Value = 435.7

# --- sheet index 4 (COM item 5) : new column 73 ---
$ws = $wb.Worksheets.Item(5)
$ws.Columns.Item(73).ColumnWidth = $ws.Columns.Item(72).ColumnWidth
$ws.Cells.Item(1, 73).Value = 20251209
$ws.Cells.Item(1, 72).Copy()
$ws.Cells.Item(1, 73).PasteSpecial(-4122)
$ws.Cells.Item(2, 73).Value = 277.18
$ws.Cells.Item(3, 73).Value = 227.92
$ws.Cells.Item(4, 73).Value = 317.08
$ws.Cells.Item(5, 73).Value = 492.02
$ws.Cells.Item(6, 73).Value = 96.71
$ws.Cells.Item(7, 73).Value = 184.97
$ws.Cells.Item(8, 73).Value = 181.84
$ws.Cells.Item(9, 73).Value = 82.28
$ws.Cells.Item(10, 73).Value = 445.17

# --- sheet index 5 (COM item 6) : new column 73 ---
$ws = $wb.Worksheets.Item(6)
$ws.Columns.Item(73).ColumnWidth = $ws.Columns.Item(72).ColumnWidth
$ws.Cells.Item(1, 73).Value = 20251209
$ws.Cells.Item(1, 72).Copy()
$ws.Cells.Item(1, 73).PasteSpecial(-4122)
$ws.Cells.Item(2, 73).Value = 32193256
$ws.Cells.Item(3, 73).Value = 25841743
$ws.Cells.Item(4, 73).Value = 30194027
$ws.Cells.Item(5, 73).Value = 14696078
$ws.Cells.Item(6, 73).Value = 51745577
$ws.Cells.Item(7, 73).Value = 144719705
$ws.Cells.Item(8, 73).Value = 18937574
$ws.Cells.Item(9, 73).Value = 10322876
$ws.Cells.Item(10, 73).Value = 62367442

# --- sheet index 6 (COM item 7) : new column 54 ---
$ws = $wb.Worksheets.Item(7)
$ws.Columns.Item(54).ColumnWidth = $ws.Columns.Item(53).ColumnWidth
$ws.Cells.Item(1, 54).Value = 20251209
$ws.Cells.Item(1, 53).Copy()
$ws.Cells.Item(1, 54).PasteSpecial(-4122)
$ws.Cells.Item(2, 54).Value = 55
$ws.Cells.Item(3, 54).Value = 34
$ws.Cells.Item(4, 54).Value = 86
$ws.Cells.Item(5, 54).Value = 51
$ws.Cells.Item(6, 54).Value = 0
$ws.Cells.Item(7, 54).Value = 47
$ws.Cells.Item(8, 54).Value = 75
$ws.Cells.Item(9, 54).Value = 0
$ws.Cells.Item(10, 54).Value = 85

# --- sheet index 7 (COM item 8) : new column 14 ---
$ws = $wb.Worksheets.Item(8)
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth
$ws.Cells.Item(1, 14).Value = 20251209
$ws.Cells.Item(1, 13).Copy()
$ws.Cells.Item(1, 14).PasteSpecial(-4122)
$ws.Cells.Item(2, 14).Value = 81
$ws.Cells.Item(3, 14).Value = 36
$ws.Cells.Item(4, 14).Value = 93
$ws.Cells.Item(5, 14).Value = 28
$ws.Cells.Item(6, 14).Value = 0
$ws.Cells.Item(7, 14).Value = 40
$ws.Cells.Item(8, 14).Value = 52
$ws.Cells.Item(9, 14).Value = 43
$ws.Cells.Item(10, 14).Value = 70

# --- sheet index 8 (COM item 9) : new column 54 ---
$ws = $wb.Worksheets.Item(9)
$ws.Columns.Item(54).ColumnWidth = $ws.Columns.Item(53).ColumnWidth
$ws.Cells.Item(1, 54).Value = 20251209
$ws.Cells.Item(1, 53).Copy()
$ws.Cells.Item(1, 54).PasteSpecial(-4122)
$ws.Cells.Item(2, 54).Value = 14
$ws.Cells.Item(3, 54).Value = -18
$ws.Cells.Item(4, 54).Value = 38
$ws.Cells.Item(5, 54).Value = 8
$ws.Cells.Item(6, 54).Value = -25
$ws.Cells.Item(7, 54).Value = 14
$ws.Cells.Item(8, 54).Value = 54
$ws.Cells.Item(9, 54).Value = -86
$ws.Cells.Item(10, 54).Value = 54

# --- sheet index 9 (COM item 10) : new column 14 ---
$ws = $wb.Worksheets.Item(10)
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth
$ws.Cells.Item(1, 14).Value = 20251209
$ws.Cells.Item(1, 13).Copy()
$ws.Cells.Item(1, 14).PasteSpecial(-4122)
$ws.Cells.Item(2, 14).Value = 57
$ws.Cells.Item(3, 14).Value = -1
$ws.Cells.Item(4, 14).Value = 81
$ws.Cells.Item(5, 14).Value = -48
$ws.Cells.Item(6, 14).Value = -83
$ws.Cells.Item(7, 14).Value = -4
$ws.Cells.Item(8, 14).Value = 14
$ws.Cells.Item(9, 14).Value = -43
$ws.Cells.Item(10, 14).Value = 31

# --- sheet index 10 (COM item 11) : new column 54 ---
$ws = $wb.Worksheets.Item(11)
$ws.Columns.Item(54).ColumnWidth = $ws.Columns.Item(53).ColumnWidth
$ws.Cells.Item(1, 54).NumberFormat = "@"
$ws.Cells.Item(1, 54).Value = "20251209"
$ws.Cells.Item(1, 53).Copy()
$ws.Cells.Item(1, 54).PasteSpecial(-4122)
$ws.Cells.Item(2, 54).Value = 101
$ws.Cells.Item(3, 54).Value = 99
$ws.Cells.Item(4, 54).Value = 104
$ws.Cells.Item(5, 54).Value = 100
$ws.Cells.Item(6, 54).Value = 31
$ws.Cells.Item(7, 54).Value = 101
$ws.Cells.Item(8, 54).Value = 106
$ws.Cells.Item(9, 54).Value = 97
$ws.Cells.Item(10, 54).Value = 105

# --- sheet index 11 (COM item 12) : new column 35 ---
$ws = $wb.Worksheets.Item(12)
$ws.Columns.Item(35).ColumnWidth = $ws.Columns.Item(34).ColumnWidth
$ws.Cells.Item(1, 35).NumberFormat = "@"
$ws.Cells.Item(1, 35).Value = "20251209"
$ws.Cells.Item(1, 34).Copy()
$ws.Cells.Item(1, 35).PasteSpecial(-4122)
$ws.Cells.Item(2, 35).Value = 20.63
$ws.Cells.Item(3, 35).Value = -32.12
$ws.Cells.Item(4, 35).Value = 20.73
$ws.Cells.Item(5, 35).Value = -12.97
$ws.Cells.Item(6, 35).Value = 15.64
$ws.Cells.Item(7, 35).Value = -38.91
$ws.Cells.Item(8, 35).Value = -19
$ws.Cells.Item(9, 35).Value = -20.19
$ws.Cells.Item(10, 35).Value = -1.02

# --- sheet index 12 (COM item 13) : new column 14 ---
$ws = $wb.Worksheets.Item(13)
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth
$ws.Cells.Item(1, 14).NumberFormat = "@"
$ws.Cells.Item(1, 14).Value = "20251209"
$ws.Cells.Item(1, 13).Copy()
$ws.Cells.Item(1, 14).PasteSpecial(-4122)
$ws.Cells.Item(2, 14).Value = 32
$ws.Cells.Item(3, 14).Value = 26
$ws.Cells.Item(4, 14).Value = 41
$ws.Cells.Item(5, 14).Value = 33
$ws.Cells.Item(6, 14).Value = 166
$ws.Cells.Item(7, 14).Value = 38
$ws.Cells.Item(8, 14).Value = 18
$ws.Cells.Item(9, 14).Value = 52
$ws.Cells.Item(10, 14).Value = 36
